$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.318.05"
$ws.Range("E2").Value = "  +3.99%  "
$ws.Range("D3").Value = "1.732.74"
$ws.Range("E3").Value = "  +2.68%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.51"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.523"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.17"
$ws.Range("E8").Value = "  +4.99%  "
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0894"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "1.977.13"
$ws.Range("E12").Value = "  +2.66%  "
$ws.Range("D13").Value = "1.730.51"
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.564"
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.74"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "28.276.66"
$ws.Range("E17").Value = "  +3.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.38"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").Value = "0.0₃0756"
$ws.Range("E19").Value = "  +1.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.95"
$ws.Range("E20").Value = "  -2.68%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.78"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.51"
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.54"
$ws.Range("E26").Value = "  +3.07%  "
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0516"
$ws.Range("E30").Value = "  +2.52%  "
$ws.Range("E31").Value = "  +2.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.43"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("D34").Value = "1.486.69"
$ws.Range("E34").Value = "  -5.88%  "
$ws.Range("E35").Value = "  -2.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.983"
$ws.Range("E36").Value = "  +3.13%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "70.39"
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.880.91"
$ws.Range("E44").Value = "  +2.17%  "
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.30"
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.797"
$ws.Range("E47").Value = "  +7.69%  "
$ws.Range("E48").Value = "  +6.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "90.82"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("E51").Value = "  -0.33%  "
